$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily-push row (row 47) under the existing data table
# (A:日付 / B:曜日 / C:時刻 / D:ランキング).
#
# Column A holds a date-like string ("2025/10/01") that must stay literal
# text (matching every other row in the sheet), not get auto-converted to
# a date serial by Excel's input heuristics. Prefixing with an apostrophe
# forces text entry, then ClearFormats() drops the resulting quote-prefix
# styling so the cell ends up with the same "no special format" look as
# the rest of the table.
$ws.Range("A47").Value = "'2025/10/01"
$ws.Range("A47").ClearFormats()

$ws.Range("B47").Value = "水"
$ws.Range("C47").Value = 20
$ws.Range("D47").Value = 197
